$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text string without altering its style,
# even when the text looks like a number (e.g. "0.1430" must keep the
# trailing zero instead of being parsed into the double 0.143).
function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '246.71'
Set-TextCell $ws.Range("D3") '21.81'
Set-TextCell $ws.Range("D4") '5.477'
Set-TextCell $ws.Range("D5") '0.05655'
Set-TextCell $ws.Range("D6") '3.378'
Set-TextCell $ws.Range("D7") '6.433'
Set-TextCell $ws.Range("D9") '1.032'
Set-TextCell $ws.Range("D10") '0.1430'
Set-TextCell $ws.Range("D11") '0.07229'
Set-TextCell $ws.Range("D12") '0.03148'
Set-TextCell $ws.Range("D13") '0.02951'
Set-TextCell $ws.Range("D14") '0.09277'
Set-TextCell $ws.Range("D15") '0.001645'
Set-TextCell $ws.Range("D16") '3.216'
Set-TextCell $ws.Range("D17") '0.04732'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell $ws.Range("D18") '0.0005853'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws.Range("D19") '0.006381'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell $ws.Range("D20") '0.005014'
$ws.Range("E20").Value = '19HotbitTokenHTBBestin24h'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell $ws.Range("D21") '0.001048'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextCell $ws.Range("D22") '0.0001502'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'UpBots'
$ws.Range("C23").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextCell $ws.Range("D23") '0.0003201'
$ws.Range("E23").Value = '22UpBotsUBXT'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws.Range("D24") '3.910'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell $ws.Range("D25") '2.109'
$ws.Range("E25").Value = '24BTSETokenBTSE'
Set-TextCell $ws.Range("D40") '0.04088'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell $ws.Range("D41") '0.1041'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell $ws.Range("D42") '0.002975'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell $ws.Range("D43") '0.003261'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
Set-TextCell $ws.Range("D44") '0.009079'
Set-TextCell $ws.Range("D45") '0.00005818'
Set-TextCell $ws.Range("D46") '0.00000000750'
Set-TextCell $ws.Range("D47") '0.7856'
Set-TextCell $ws.Range("D48") '0.01602'
$ws.Range("E48").Value = '47BOLOBOLO'
Set-TextCell $ws.Range("D49") '0.00002101'
Set-TextCell $ws.Range("D50") '0.01010'
